# Populate the AprilRaw sheet with the raw monthly statistics data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AprilRaw")

$ws.Range("A1").Value = 'Library'
$ws.Range("B1").Value = 'Items owned by this library checked out at this library this month'
$ws.Range("C1").Value = 'Items owned by other libraries checked out at this library this month'
$ws.Range("D1").Value = 'Total circulation this month'

$ws.Range("A2").Value = 'Atchison Public Library'
$ws.Range("B2").Value = 4060
$ws.Range("C2").Value = 1339
$ws.Range("D2").Value = 5399

$ws.Range("A3").Value = 'Baldwin City Public Library'
$ws.Range("B3").Value = 2703
$ws.Range("C3").Value = 603
$ws.Range("D3").Value = 3306

$ws.Range("A4").Value = 'Basehor Community Library'
$ws.Range("B4").Value = 7973
$ws.Range("C4").Value = 1106
$ws.Range("D4").Value = 9079

$ws.Range("A5").Value = 'Bern Community Library'
$ws.Range("B5").Value = 102
$ws.Range("C5").Value = 40
$ws.Range("D5").Value = 142

$ws.Range("A6").Value = 'Bonner Springs City Library'
$ws.Range("B6").Value = 4794
$ws.Range("C6").Value = 991
$ws.Range("D6").Value = 5785

$ws.Range("A7").Value = 'Burlingame Community Library'
$ws.Range("B7").Value = 540
$ws.Range("C7").Value = 501
$ws.Range("D7").Value = 1041

$ws.Range("A8").Value = 'Carbondale City Library'
$ws.Range("B8").Value = 519
$ws.Range("C8").Value = 120
$ws.Range("D8").Value = 639

$ws.Range("A9").Value = 'Centralia Community Library'
$ws.Range("B9").Value = 311
$ws.Range("C9").Value = 51
$ws.Range("D9").Value = 362

$ws.Range("A10").Value = 'Corning City Library'
$ws.Range("B10").Value = 72
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 73

$ws.Range("A11").Value = 'Digital Content'

$ws.Range("A12").Value = 'Doniphan County Library - Elwood'
$ws.Range("B12").Value = 72
$ws.Range("C12").Value = 32
$ws.Range("D12").Value = 104

$ws.Range("A13").Value = 'Doniphan County Library - Highland'
$ws.Range("B13").Value = 189
$ws.Range("C13").Value = 148
$ws.Range("D13").Value = 337

$ws.Range("A14").Value = 'Doniphan County Library - Troy'
$ws.Range("B14").Value = 451
$ws.Range("C14").Value = 164
$ws.Range("D14").Value = 615

$ws.Range("A15").Value = 'Doniphan County Library - Wathena'
$ws.Range("B15").Value = 481
$ws.Range("C15").Value = 47
$ws.Range("D15").Value = 528

$ws.Range("A16").Value = 'Effingham Community Library'
$ws.Range("B16").Value = 295
$ws.Range("C16").Value = 81
$ws.Range("D16").Value = 376

$ws.Range("A17").Value = 'Eudora Community Library'
$ws.Range("B17").Value = 1654
$ws.Range("C17").Value = 514
$ws.Range("D17").Value = 2168

$ws.Range("A18").Value = 'Everest, Barnes Reading Room'
$ws.Range("B18").Value = 63
$ws.Range("C18").Value = 76
$ws.Range("D18").Value = 139

$ws.Range("A19").Value = 'Hiawatha, Morrill Public Library'
$ws.Range("B19").Value = 1728
$ws.Range("C19").Value = 637
$ws.Range("D19").Value = 2365

$ws.Range("A20").Value = 'Highland Community College'
$ws.Range("B20").Value = 31
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 32

$ws.Range("A21").Value = 'Holton, Beck-Bookman Library'
$ws.Range("B21").Value = 1531
$ws.Range("C21").Value = 410
$ws.Range("D21").Value = 1941

$ws.Range("A22").Value = 'Horton Public Library'
$ws.Range("B22").Value = 117
$ws.Range("C22").Value = 40
$ws.Range("D22").Value = 157

$ws.Range("A23").Value = 'Lansing Community Library'
$ws.Range("B23").Value = 1630
$ws.Range("C23").Value = 671
$ws.Range("D23").Value = 2301

$ws.Range("A24").Value = 'Leavenworth Public Library'
$ws.Range("B24").Value = 8588
$ws.Range("C24").Value = 1739
$ws.Range("D24").Value = 10327

$ws.Range("A25").Value = 'Linwood Community Library'
$ws.Range("B25").Value = 606
$ws.Range("C25").Value = 137
$ws.Range("D25").Value = 743

$ws.Range("A26").Value = 'Louisburg Library'

$ws.Range("A27").Value = 'Lyndon Carnegie Library'
$ws.Range("B27").Value = 439
$ws.Range("C27").Value = 223
$ws.Range("D27").Value = 662

$ws.Range("A28").Value = 'McLouth Public Library'
$ws.Range("B28").Value = 174
$ws.Range("C28").Value = 114
$ws.Range("D28").Value = 288

$ws.Range("A29").Value = 'Meriden-Ozawkie Public Library'
$ws.Range("B29").Value = 1341
$ws.Range("C29").Value = 576
$ws.Range("D29").Value = 1917

$ws.Range("A30").Value = 'Northeast Kansas Library System'
$ws.Range("B30").Value = 37
$ws.Range("C30").Value = 45
$ws.Range("D30").Value = 82

$ws.Range("A31").Value = 'Nortonville Public Library'
$ws.Range("B31").Value = 309
$ws.Range("C31").Value = 82
$ws.Range("D31").Value = 391

$ws.Range("A32").Value = 'Osage City Library'
$ws.Range("B32").Value = 1279
$ws.Range("C32").Value = 460
$ws.Range("D32").Value = 1739

$ws.Range("A33").Value = 'Osawatomie Public Library'
$ws.Range("B33").Value = 1004
$ws.Range("C33").Value = 300
$ws.Range("D33").Value = 1304

$ws.Range("A34").Value = 'Oskaloosa Public Library'
$ws.Range("B34").Value = 475
$ws.Range("C34").Value = 154
$ws.Range("D34").Value = 629

$ws.Range("A35").Value = 'Ottawa Library'
$ws.Range("B35").Value = 6234
$ws.Range("C35").Value = 831
$ws.Range("D35").Value = 7065

$ws.Range("A36").Value = 'Overbrook Public Library'
$ws.Range("B36").Value = 734
$ws.Range("C36").Value = 178
$ws.Range("D36").Value = 912

$ws.Range("A37").Value = 'Paola Free Library'
$ws.Range("B37").Value = 2970
$ws.Range("C37").Value = 428
$ws.Range("D37").Value = 3398

$ws.Range("A38").Value = 'Perry-Lecompton Community Library'
$ws.Range("B38").Value = 63
$ws.Range("C38").Value = 35
$ws.Range("D38").Value = 98

$ws.Range("A39").Value = 'Pomona Community Library'
$ws.Range("B39").Value = 40
$ws.Range("C39").Value = 70
$ws.Range("D39").Value = 110

$ws.Range("A40").Value = 'Prairie Hills Schools - Axtell Public School'
$ws.Range("B40").Value = 641
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 643

$ws.Range("A41").Value = 'Prairie Hills Schools - Sabetha Elementary School'
$ws.Range("B41").Value = 1763
$ws.Range("C41").Value = 76
$ws.Range("D41").Value = 1839

$ws.Range("A42").Value = 'Prairie Hills Schools - Sabetha High School'
$ws.Range("B42").Value = 27
$ws.Range("C42").Value = 9
$ws.Range("D42").Value = 36

$ws.Range("A43").Value = 'Prairie Hills Schools - Sabetha Middle School'
$ws.Range("B43").Value = 118
$ws.Range("C43").Value = 10
$ws.Range("D43").Value = 128

$ws.Range("A44").Value = 'Prairie Hills Schools - Wetmore Academic Center (Permanently closed)'

$ws.Range("A45").Value = 'Richmond Public Library'
$ws.Range("B45").Value = 272
$ws.Range("C45").Value = 63
$ws.Range("D45").Value = 335

$ws.Range("A46").Value = 'Rossville Community Library'
$ws.Range("B46").Value = 1260
$ws.Range("C46").Value = 477
$ws.Range("D46").Value = 1737

$ws.Range("A47").Value = 'Sabetha, Mary Cotton Library'
$ws.Range("B47").Value = 3132
$ws.Range("C47").Value = 1104
$ws.Range("D47").Value = 4236

$ws.Range("A48").Value = 'Seneca Free Library'
$ws.Range("B48").Value = 1450
$ws.Range("C48").Value = 273
$ws.Range("D48").Value = 1723

$ws.Range("A49").Value = 'Silver Lake Library'
$ws.Range("B49").Value = 1049
$ws.Range("C49").Value = 528
$ws.Range("D49").Value = 1577

$ws.Range("A50").Value = 'Tonganoxie Public Library'
$ws.Range("B50").Value = 2981
$ws.Range("C50").Value = 1008
$ws.Range("D50").Value = 3989

$ws.Range("A51").Value = 'Valley Falls, Delaware Township Library'
$ws.Range("B51").Value = 421
$ws.Range("C51").Value = 214
$ws.Range("D51").Value = 635

$ws.Range("A52").Value = 'Wellsville City Library'
$ws.Range("B52").Value = 898
$ws.Range("C52").Value = 295
$ws.Range("D52").Value = 1193

$ws.Range("A53").Value = 'Wetmore Public Library'
$ws.Range("B53").Value = 134
$ws.Range("C53").Value = 117
$ws.Range("D53").Value = 251

$ws.Range("A54").Value = 'Williamsburg Community Library'
$ws.Range("B54").Value = 331
$ws.Range("C54").Value = 28
$ws.Range("D54").Value = 359

$ws.Range("A55").Value = 'Winchester Public Library'
$ws.Range("B55").Value = 215
$ws.Range("C55").Value = 275
$ws.Range("D55").Value = 490

# Select the April sheet as the active tab (matches tabSelected/activeTab move).
$april = $wb.Worksheets.Item("April")
$april.Activate()

